$d = $word.ActiveDocument

$replacements = @(
    @{old="2026-01-30 Friday"; new="2026-01-31 Saturday"},
    @{old="32×88="; new="56×96="},
    @{old="52×35="; new="31×86="},
    @{old="97×76="; new="13×28="},
    @{old="31×55="; new="94×71="},
    @{old="52×78="; new="52×50="},
    @{old="79×95="; new="34×66="},
    @{old="33×51="; new="67×40="},
    @{old="53×52="; new="46×90="},
    @{old="91×47="; new="95×54="},
    @{old="46×99="; new="67×37="},
    @{old="75×86="; new="89×98="},
    @{old="95×33="; new="35×35="},
    @{old="68×82="; new="99×34="},
    @{old="92×66="; new="80×80="},
    @{old="26×37="; new="69×84="},
    @{old="79×55="; new="55×49="},
    @{old="92×32="; new="99×59="},
    @{old="22×95="; new="99×51="},
    @{old="40×83="; new="45×24="},
    @{old="67×45="; new="87×40="},
    @{old="35×43="; new="82×41="},
    @{old="42×17="; new="40×50="},
    @{old="26×48="; new="64×14="},
    @{old="47×71="; new="50×11="},
    @{old="83×73="; new="41×11="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
